$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices, volumes, and the row-45/46 coin swap).
# Column D values are prefixed with a leading apostrophe so Excel keeps them as
# literal text (matching the source inlineStr cells) instead of re-parsing them
# as numbers and losing formatting such as "1.000", "6.550", "29.638.54", etc.

# Row 2
$ws.Range("D2").Value = "'29.638.54"
$ws.Range("E2").Value = "  +2.62%  "
# Row 3
$ws.Range("D3").Value = "'1.860.90"
$ws.Range("E3").Value = "  +1.96%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("D5").Value = "'245.38"
$ws.Range("E5").Value = "  +2.24%  "
# Row 6
$ws.Range("D6").Value = "'0.6975"
$ws.Range("E6").Value = "  +1.12%  "
# Row 7
$ws.Range("E7").Value = "  +0.04%  "
# Row 8
$ws.Range("D8").Value = "'0.07726"
$ws.Range("E8").Value = "  +1.64%  "
# Row 9
$ws.Range("D9").Value = "'0.3066"
$ws.Range("E9").Value = "  +1.67%  "
# Row 10
$ws.Range("D10").Value = "'23.67"
$ws.Range("E10").Value = "  +1.25%  "
# Row 11
$ws.Range("D11").Value = "'0.07765"
$ws.Range("E11").Value = "  +0.38%  "
# Row 12
$ws.Range("D12").Value = "'5.164"
$ws.Range("E12").Value = "  +2.53%  "
# Row 13
$ws.Range("D13").Value = "'1.856.08"
$ws.Range("E13").Value = "  +1.84%  "
# Row 14
$ws.Range("D14").Value = "'92.31"
$ws.Range("E14").Value = "  +2.55%  "
# Row 15
$ws.Range("D15").Value = "'0.6932"
$ws.Range("E15").Value = "  +3.31%  "
# Row 16
$ws.Range("D16").Value = "'6.550"
$ws.Range("E16").Value = "  +3.14%  "
# Row 17
$ws.Range("D17").Value = "'29.626.43"
# Row 18
$ws.Range("D18").Value = "'0.000008355"
$ws.Range("E18").Value = "  +1.15%  "
# Row 19
$ws.Range("D19").Value = "'2.106.28"
$ws.Range("E19").Value = "  +1.61%  "
# Row 20
$ws.Range("D20").Value = "'241.85"
$ws.Range("E20").Value = "  -0.24%  "
# Row 21
$ws.Range("D21").Value = "'12.76"
$ws.Range("E21").Value = "  +1.35%  "
# Row 22
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "
# Row 23
$ws.Range("E23").Value = "  +3.07%  "
# Row 24
$ws.Range("E24").Value = "  +0.04%  "
# Row 25
$ws.Range("D25").Value = "'0.1507"
$ws.Range("E25").Value = "  +2.58%  "
# Row 26
$ws.Range("D26").Value = "'8.908"
$ws.Range("E26").Value = "  +2.36%  "
# Row 27
$ws.Range("D27").Value = "'159.65"
$ws.Range("E27").Value = "  -0.68%  "
# Row 28
$ws.Range("D28").Value = "'18.30"
$ws.Range("E28").Value = "  +0.91%  "
# Row 29
$ws.Range("D29").Value = "'1.535"
$ws.Range("E29").Value = "  +0.31%  "
# Row 30
$ws.Range("D30").Value = "'4.253"
$ws.Range("E30").Value = "  +1.61%  "
# Row 31
$ws.Range("D31").Value = "'4.187"
$ws.Range("E31").Value = "  +1.78%  "
# Row 32
$ws.Range("E32").Value = "  +0.31%  "
# Row 33
$ws.Range("D33").Value = "'0.05102"
$ws.Range("E33").Value = "  +0.34%  "
# Row 34
$ws.Range("D34").Value = "'0.7808"
$ws.Range("E34").Value = "  +3.96%  "
# Row 35
$ws.Range("D35").Value = "'1.901"
$ws.Range("E35").Value = "  +5.38%  "
# Row 36
$ws.Range("D36").Value = "'1.156"
$ws.Range("E36").Value = "  +1.77%  "
# Row 37
$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "  +0.15%  "
# Row 38
$ws.Range("D38").Value = "'1.325.59"
$ws.Range("E38").Value = "  +10.45%  "
# Row 39
$ws.Range("D39").Value = "'0.01877"
$ws.Range("E39").Value = "  +2.36%  "
# Row 40
$ws.Range("D40").Value = "'2.732"
$ws.Range("E40").Value = "  +1.97%  "
# Row 41
$ws.Range("D41").Value = "'0.9610"
$ws.Range("E41").Value = "  +3.98%  "
# Row 42
$ws.Range("D42").Value = "'106.48"
$ws.Range("E42").Value = "  -1.33%  "
# Row 43
$ws.Range("E43").Value = "  +13.76%  "
# Row 44
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  +0.05%  "
# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000126"
$ws.Range("E45").Value = "  +4.13%  "
# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'9.778"
$ws.Range("E46").Value = "  +3.44%  "
# Row 47
$ws.Range("D47").Value = "'2.005.27"
$ws.Range("E47").Value = "  +1.49%  "
# Row 48
$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = "  +1.04%  "
# Row 49
$ws.Range("D49").Value = "'1.787"
$ws.Range("E49").Value = "  +4.03%  "
# Row 50
$ws.Range("D50").Value = "'64.51"
$ws.Range("E50").Value = "  +4.03%  "
# Row 51
$ws.Range("D51").Value = "'6.995"
$ws.Range("E51").Value = "  +2.01%  "
